# Update "想去人数" (F column) values across the four sheets to match the
# refreshed scrape output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Column F is the 6th column on every sheet in this workbook.
$col = 6

# 展览 (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    3  = 1014
    5  = 11
    6  = 456
    7  = 732
    11 = 406
    12 = 212
    14 = 849
    16 = 1994
    17 = 491
    18 = 7421
    19 = 547
    21 = 59
    22 = 94
    24 = 223
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, $col).Value = $sheet1Updates[$row]
}

# 演出 (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$sheet2Updates = @{
    2  = 521
    8  = 121
    10 = 4
}
foreach ($row in $sheet2Updates.Keys) {
    $ws2.Cells.Item($row, $col).Value = $sheet2Updates[$row]
}

# 本地生活 (sheet3)
$ws3 = $wb.Worksheets.Item("本地生活")
$sheet3Updates = @{
    2 = 5514
    3 = 399
    4 = 388
}
foreach ($row in $sheet3Updates.Keys) {
    $ws3.Cells.Item($row, $col).Value = $sheet3Updates[$row]
}

# 全部类型 (sheet4) - combined listing of all the above
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    3  = 5514
    4  = 399
    5  = 388
    6  = 521
    7  = 1014
    11 = 11
    12 = 456
    13 = 732
    18 = 406
    19 = 212
    23 = 849
    25 = 121
    26 = 1994
    27 = 491
    28 = 7421
    30 = 4
    31 = 547
    33 = 59
    34 = 94
    37 = 223
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, $col).Value = $sheet4Updates[$row]
}
